$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2454.1667
$ws.Range("I132").Value = 2562.1052
$ws.Range("J132").Value = 2044
$ws.Range("K132").Value = 7686.3156
$ws.Range("L132").Value = 6132
$ws.Range("M132").Value = -5156.3156
$ws.Range("N132").Value = -11192
$ws.Range("H137").Value = 1377.0769
$ws.Range("I137").Value = 1000.3333
$ws.Range("J137").Value = 2224.75
$ws.Range("K137").Value = 3000.9999
$ws.Range("L137").Value = 6674.25
$ws.Range("M137").Value = -450.9998999999998
$ws.Range("N137").Value = -11774.25

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1611.5333
$ws.Range("I2").Value = 861.1579
$ws.Range("K2").Value = 861.1579
$ws.Range("M2").Value = -748.1579
$ws.Range("H37").Value = 7000.5
$ws.Range("I37").Value = 3400
$ws.Range("J37").Value = 8200.666999999999
$ws.Range("K37").Value = 3400
$ws.Range("L37").Value = 8200.666999999999
$ws.Range("M37").Value = -3127
$ws.Range("N37").Value = -8746.666999999999
$ws.Range("H74").Value = 4494.1113
$ws.Range("I74").Value = 5064.5713
$ws.Range("J74").Value = 2497.5
$ws.Range("K74").Value = 5064.5713
$ws.Range("L74").Value = 2497.5
$ws.Range("M74").Value = -4190.5713
$ws.Range("N74").Value = -4245.5
$ws.Range("H77").Value = 4494.1113
$ws.Range("I77").Value = 5064.5713
$ws.Range("J77").Value = 2497.5
$ws.Range("K77").Value = 25322.8565
$ws.Range("L77").Value = 12487.5
$ws.Range("M77").Value = -20954.8565
$ws.Range("N77").Value = -21223.5
$ws.Range("H116").Value = 1611.5333
$ws.Range("I116").Value = 861.1579
$ws.Range("K116").Value = 861.1579
$ws.Range("M116").Value = 1432.8421

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1611.5333
$ws.Range("I3").Value = 861.1579
$ws.Range("K3").Value = 861.1579
$ws.Range("M3").Value = -747.1579

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 8815.666999999999
$ws.Range("J50").Value = 8815.666999999999
$ws.Range("L50").Value = 8815.666999999999
$ws.Range("N50").Value = -10065.667
$ws.Range("H51").Value = 9500.75
$ws.Range("J51").Value = 9500.75
$ws.Range("L51").Value = 9500.75
$ws.Range("N51").Value = -10972.75
$ws.Range("H58").Value = 920.2222
$ws.Range("I58").Value = 644.9286
$ws.Range("J58").Value = 1216.6923
$ws.Range("K58").Value = 644.9286
$ws.Range("L58").Value = 1216.6923
$ws.Range("M58").Value = -441.9286
$ws.Range("N58").Value = -1622.6923
$ws.Range("H59").Value = 12518.333
$ws.Range("J59").Value = 12518.333
$ws.Range("L59").Value = 12518.333
$ws.Range("N59").Value = -14808.333
$ws.Range("H61").Value = 9500.75
$ws.Range("J61").Value = 9500.75
$ws.Range("L61").Value = 9500.75
$ws.Range("N61").Value = -10196.75
$ws.Range("H68").Value = 16599.8
$ws.Range("J68").Value = 16599.8
$ws.Range("L68").Value = 16599.8
$ws.Range("N68").Value = -18097.8
$ws.Range("H71").Value = 16599.8
$ws.Range("J71").Value = 16599.8
$ws.Range("L71").Value = 49799.39999999999
$ws.Range("N71").Value = -57287.39999999999
$ws.Range("H74").Value = 15793.444
$ws.Range("I74").Value = 2285
$ws.Range("J74").Value = 17482
$ws.Range("K74").Value = 2285
$ws.Range("L74").Value = 17482
$ws.Range("N74").Value = -19230
$ws.Range("M74").Value = -1411
$ws.Range("H77").Value = 15793.444
$ws.Range("I77").Value = 2285
$ws.Range("J77").Value = 17482
$ws.Range("K77").Value = 6855
$ws.Range("L77").Value = 52446
$ws.Range("N77").Value = -61182
$ws.Range("M77").Value = -2487
$ws.Range("H80").Value = 24633.334
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("H83").Value = 24633.334
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("H134").Value = 1099.3158
$ws.Range("I134").Value = 1077.0555
$ws.Range("K134").Value = 3231.1665
$ws.Range("M134").Value = -696.1664999999998
$ws.Range("H136").Value = 920.2222
$ws.Range("I136").Value = 644.9286
$ws.Range("J136").Value = 1216.6923
$ws.Range("K136").Value = 1934.7858
$ws.Range("L136").Value = 3650.0769
$ws.Range("M136").Value = 615.2142000000001
$ws.Range("N136").Value = -8750.0769
$ws.Range("M80").ClearContents()
$ws.Range("M83").ClearContents()

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 972.0625
$ws.Range("I5").Value = 809.36365
$ws.Range("J5").Value = 1330
$ws.Range("K5").Value = 2428.09095
$ws.Range("L5").Value = 3990
$ws.Range("M5").Value = -2316.09095
$ws.Range("N5").Value = -4214
$ws.Range("H124").Value = 4477.778
$ws.Range("I124").Value = 800
$ws.Range("J124").Value = 4937.5
$ws.Range("K124").Value = 2400
$ws.Range("L124").Value = 14812.5
$ws.Range("N124").Value = -24632.5
$ws.Range("M124").Value = 2510
$ws.Range("H131").Value = 5747985
$ws.Range("I131").Value = 976.8125
$ws.Range("J131").Value = 7043085.5
$ws.Range("K131").Value = 2930.4375
$ws.Range("L131").Value = 21129256.5
$ws.Range("M131").Value = 2109.5625
$ws.Range("N131").Value = -21139336.5
$ws.Range("H135").Value = 972.0625
$ws.Range("I135").Value = 809.36365
$ws.Range("J135").Value = 1330
$ws.Range("K135").Value = 7284.27285
$ws.Range("L135").Value = 11970
$ws.Range("M135").Value = -4749.27285
$ws.Range("N135").Value = -17040

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 36653.324
$ws.Range("I132").Value = 44438.2
$ws.Range("J132").Value = 4216.3335
$ws.Range("K132").Value = 133314.6
$ws.Range("L132").Value = 12649.0005
$ws.Range("M132").Value = -130784.6
$ws.Range("N132").Value = -17709.0005
$ws.Range("H140").Value = 37088
$ws.Range("J140").Value = 37088
$ws.Range("L140").Value = 37088
$ws.Range("N140").Value = -47448

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3054.111
$ws.Range("I40").Value = 2125
$ws.Range("K40").Value = 2125
$ws.Range("M40").Value = -1989
$ws.Range("H122").Value = 2384.5
$ws.Range("I122").Value = 2365.4285
$ws.Range("J122").Value = 2429
$ws.Range("K122").Value = 7096.2855
$ws.Range("L122").Value = 7287
$ws.Range("M122").Value = -4646.2855
$ws.Range("N122").Value = -12187

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 3384.4
$ws.Range("I37").Value = 2526
$ws.Range("K37").Value = 2526
$ws.Range("M37").Value = -2323
$ws.Range("H122").Value = 1608.6666
$ws.Range("I122").Value = 1563
$ws.Range("J122").Value = 1700
$ws.Range("K122").Value = 4689
$ws.Range("L122").Value = 5100
$ws.Range("M122").Value = -2239
$ws.Range("N122").Value = -10000
$ws.Range("H132").Value = 20002322
$ws.Range("I132").Value = 23438242
$ws.Range("J132").Value = 11518
$ws.Range("K132").Value = 70314726
$ws.Range("L132").Value = 34554
$ws.Range("M132").Value = -70312196
$ws.Range("N132").Value = -39614
$ws.Range("H136").Value = 1050.129
$ws.Range("I136").Value = 746.04346
$ws.Range("J136").Value = 1924.375
$ws.Range("K136").Value = 2238.13038
$ws.Range("L136").Value = 5773.125
$ws.Range("M136").Value = 311.8696199999999
$ws.Range("N136").Value = -10873.125
